# Upload new version with timestamp
# Populate the transactions table (rows 4..17) with the new report data,
# push the existing summary row (old row 5) and footer row (old row 6)
# down to rows 18 and 19, and fill in the corresponding values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: insert 13 new rows before the old "summary" row (row 5)
#    so the single empty data row (row 4) becomes 14 data rows (4..17).
#    The old summary row (was 5) and footer row (was 6) shift to 18/19.
# ---------------------------------------------------------------------
$ws.Range("A5:A17").EntireRow.Insert()

# Give every new row (5..17) the same formatting as the template row (4):
# number format, font, fill, border, alignment.
$ws.Range("A4:N4").Copy()
for ($r = 5; $r -le 17; $r++) {
    $ws.Range("A" + $r + ":N" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Re-create the merges lost on the unmerged freshly-inserted rows.
for ($r = 4; $r -le 17; $r++) {
    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
}

# ---------------------------------------------------------------------
# 2) Row heights: rows alternate 24.75 / 25.5 as in the source rows,
#    the summary row is 25.5 and the footer row is 17.25.
# ---------------------------------------------------------------------
$rowHeights = @{
    4 = 24.75; 5 = 25.5; 6 = 24.75; 7 = 25.5; 8 = 25.5; 9 = 24.75;
    10 = 25.5; 11 = 24.75; 12 = 25.5; 13 = 25.5; 14 = 24.75; 15 = 25.5;
    16 = 24.75; 17 = 25.5; 18 = 25.5; 19 = 17.25
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# Columns B (name) and H (balance) hold strings like "0:0" that must stay
# text, not be reinterpreted as times.
$ws.Range("B4:B17").NumberFormat = "@"
$ws.Range("H4:H17").NumberFormat = "@"

# ---------------------------------------------------------------------
# 3) Fill in the 14 transaction rows.
# ---------------------------------------------------------------------
$rows = @(
    @{ Row = 4;  Num = 1;  Name = "ANGIOFOX (EFFOX) 25MG LONG 30 CAPS.";         Balance = "0:0";    Price = 114;   Count = 1 },
    @{ Row = 5;  Num = 2;  Name = "AUGMENTIN 457MG/5ML SUSP. 70 ML";             Balance = "1:0";    Price = 137;   Count = 1 },
    @{ Row = 6;  Num = 3;  Name = "BLOKATENS 10/160MG 28 F.C.TABS.";             Balance = "0:0";    Price = 160;   Count = 1 },
    @{ Row = 7;  Num = 4;  Name = "COLOVATIL 30 F.C. TABS";                      Balance = "0:0";    Price = 63;    Count = 1 },
    @{ Row = 8;  Num = 5;  Name = "GAVISCON LIQUID 24 SACHETS 10 ML";            Balance = "0:9";    Price = 12;    Count = 0.04 },
    @{ Row = 9;  Num = 6;  Name = "GINKGO BILOBA 30 CAPS.";                      Balance = "0:0";    Price = 186;   Count = 1 },
    @{ Row = 10; Num = 7;  Name = "MILGA ADVANCE 30 F.C. TABS";                  Balance = "0:0";    Price = 136.5; Count = 1 },
    @{ Row = 11; Num = 8;  Name = "PERLOC 40MG 14 F.C.TAB.";                     Balance = "0:0";    Price = 68.25; Count = 1 },
    @{ Row = 12; Num = 9;  Name = "RHINEX 0.05% INFANTILE NASAL DROPS 10 ML";    Balance = "2:0";    Price = 18;    Count = 1 },
    @{ Row = 13; Num = 10; Name = "RIVO 320MG 20*10 TABS";                       Balance = "1:2";    Price = 14.1;  Count = 0.1 },
    @{ Row = 14; Num = 11; Name = "VASTAREL MR 35MG 30 F.C.TAB.";                Balance = "2:0";    Price = 175;   Count = 1 },
    @{ Row = 15; Num = 12; Name = "WATER FOR INJECTION AMP. 5 ML";               Balance = "7816:0"; Price = 2.5;   Count = 1 },
    @{ Row = 16; Num = 13; Name = "سويت كوكو";                                   Balance = "22:0";   Price = 25;    Count = 1 },
    @{ Row = 17; Num = 14; Name = "مرطب شفاه لونا جوز هند ابيض";                 Balance = "3:0";    Price = 20;    Count = 1 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A" + $r).Value = $item.Num
    $ws.Range("B" + $r).Value = $item.Name
    $ws.Range("H" + $r).Value = $item.Balance
    $ws.Range("L" + $r).Value = $item.Price
    $ws.Range("N" + $r).Value = $item.Count
}

# ---------------------------------------------------------------------
# 4) Totals row (was row 5, now row 18): sum of the price column.
# ---------------------------------------------------------------------
$ws.Range("K18").Value = 1131.35

Write-Output "done"
